$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 33433.43
$ws.Range("I40").Value = 27758.5
$ws.Range("J40").Value = 41000
$ws.Range("K40").Value = 27758.5
$ws.Range("L40").Value = 41000
$ws.Range("M40").Value = -27583.5
$ws.Range("N40").Value = -41350
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = $null
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = $null
$ws.Range("H121").Value = 4577
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 4577
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 13731
$ws.Range("N121").Value = -17225
$ws.Range("H132").Value = 107708.69
$ws.Range("I132").Value = 275567.88
$ws.Range("J132").Value = 13821.339
$ws.Range("K132").Value = 826703.64
$ws.Range("L132").Value = 41464.017
$ws.Range("M132").Value = -824173.64
$ws.Range("N132").Value = -46524.017
$ws.Range("H139").Value = 99640
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 99640
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 99640
$ws.Range("N139").Value = -109920
$ws.Range("H140").Value = 68051
$ws.Range("I140").Value = 70000
$ws.Range("J140").Value = 67834.44500000001
$ws.Range("K140").Value = 70000
$ws.Range("L140").Value = 67834.44500000001
$ws.Range("M140").Value = -64820
$ws.Range("N140").Value = -78194.44500000001
$ws.Range("H141").Value = 4891.6924
$ws.Range("I141").Value = 4417.4546
$ws.Range("J141").Value = 7500
$ws.Range("K141").Value = 13252.3638
$ws.Range("L141").Value = 22500
$ws.Range("M141").Value = -8072.363799999999
$ws.Range("N141").Value = -32860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 987.46155
$ws.Range("I102").Value = 1007.75
$ws.Range("J102").Value = 744
$ws.Range("K102").Value = 1007.75
$ws.Range("L102").Value = 744
$ws.Range("M102").Value = 614.25
$ws.Range("N102").Value = -3988
$ws.Range("H110").Value = 4624.095
$ws.Range("I110").Value = 4050.9443
$ws.Range("J110").Value = 8063
$ws.Range("K110").Value = 4050.9443
$ws.Range("L110").Value = 8063
$ws.Range("M110").Value = -2005.9443
$ws.Range("N110").Value = -12153
$ws.Range("H140").Value = 112649.75
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 112649.75
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 112649.75
$ws.Range("N140").Value = -123009.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3609.2
$ws.Range("I20").Value = 3112.7273
$ws.Range("J20").Value = 4974.5
$ws.Range("K20").Value = 3112.7273
$ws.Range("L20").Value = 4974.5
$ws.Range("M20").Value = -2865.7273
$ws.Range("N20").Value = -5468.5
$ws.Range("H92").Value = 119999.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 119999.5
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 119999.5
$ws.Range("N92").Value = -124991.5
$ws.Range("H94").Value = 23810114
$ws.Range("I94").Value = 33333760
$ws.Range("J94").Value = 999.8333
$ws.Range("K94").Value = 33333760
$ws.Range("L94").Value = 999.8333
$ws.Range("M94").Value = -33333309
$ws.Range("N94").Value = -1901.8333
$ws.Range("H130").Value = 62694.5
$ws.Range("I130").Value = 62709
$ws.Range("J130").Value = 62680
$ws.Range("K130").Value = 62709
$ws.Range("L130").Value = 62680
$ws.Range("M130").Value = -57689
$ws.Range("N130").Value = -72720

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3890.676
$ws.Range("I31").Value = 2745.0576
$ws.Range("J31").Value = 7026.0527
$ws.Range("K31").Value = 2745.0576
$ws.Range("L31").Value = 7026.0527
$ws.Range("M31").Value = -2450.0576
$ws.Range("N31").Value = -7616.0527
$ws.Range("H34").Value = 3890.676
$ws.Range("I34").Value = 2745.0576
$ws.Range("J34").Value = 7026.0527
$ws.Range("K34").Value = 2745.0576
$ws.Range("L34").Value = 7026.0527
$ws.Range("M34").Value = -2543.0576
$ws.Range("N34").Value = -7430.0527

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 25562
$ws.Range("I74").Value = 20998
$ws.Range("J74").Value = 27083.334
$ws.Range("K74").Value = 62994
$ws.Range("L74").Value = 81250.00199999999
$ws.Range("M74").Value = -61933
$ws.Range("N74").Value = -83372.00199999999
$ws.Range("H77").Value = 25562
$ws.Range("I77").Value = 20998
$ws.Range("J77").Value = 27083.334
$ws.Range("K77").Value = 188982
$ws.Range("L77").Value = 243750.006
$ws.Range("M77").Value = -183678
$ws.Range("N77").Value = -254358.006
$ws.Range("H122").Value = 576.1667
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 576.1667
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 5185.5003
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = -10085.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 26000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 26000
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 26000
$ws.Range("N48").Value = -26970
$ws.Range("H70").Value = 7826.303
$ws.Range("I70").Value = 7267.8696
$ws.Range("J70").Value = 9110.700000000001
$ws.Range("K70").Value = 7267.8696
$ws.Range("L70").Value = 9110.700000000001
$ws.Range("M70").Value = -6997.8696
$ws.Range("N70").Value = -9650.700000000001
$ws.Range("H73").Value = 7826.303
$ws.Range("I73").Value = 7267.8696
$ws.Range("J73").Value = 9110.700000000001
$ws.Range("K73").Value = 7267.8696
$ws.Range("L73").Value = 9110.700000000001
$ws.Range("M73").Value = -6331.8696
$ws.Range("N73").Value = -10982.7
$ws.Range("H107").Value = 673.6
$ws.Range("I107").Value = 694
$ws.Range("J107").Value = 643
$ws.Range("K107").Value = 694
$ws.Range("L107").Value = 643
$ws.Range("M107").Value = 1226
$ws.Range("N107").Value = -4483
$ws.Range("H113").Value = 528360.8
$ws.Range("I113").Value = 1183800.6
$ws.Range("J113").Value = 4009
$ws.Range("K113").Value = 1183800.6
$ws.Range("L113").Value = 4009
$ws.Range("M113").Value = -1181630.6
$ws.Range("N113").Value = -8349
$ws.Range("H126").Value = 2372.4285
$ws.Range("I126").Value = 2190.1
$ws.Range("J126").Value = 2828.25
$ws.Range("K126").Value = 6570.299999999999
$ws.Range("L126").Value = 8484.75
$ws.Range("M126").Value = -4100.299999999999
$ws.Range("N126").Value = -13424.75
$ws.Range("H132").Value = 5896.6294
$ws.Range("I132").Value = 4991.579
$ws.Range("J132").Value = 8046.125
$ws.Range("K132").Value = 14974.737
$ws.Range("L132").Value = 24138.375
$ws.Range("M132").Value = -12444.737
$ws.Range("N132").Value = -29198.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2253.1304
$ws.Range("I61").Value = 2283.8823
$ws.Range("J61").Value = 2166
$ws.Range("K61").Value = 2283.8823
$ws.Range("L61").Value = 2166
$ws.Range("M61").Value = -2081.8823
$ws.Range("N61").Value = -2570
$ws.Range("H93").Value = 1111.8125
$ws.Range("I93").Value = 1297.7778
$ws.Range("J93").Value = 872.7143
$ws.Range("K93").Value = 1297.7778
$ws.Range("L93").Value = 872.7143
$ws.Range("M93").Value = -49.77780000000007
$ws.Range("N93").Value = -3368.7143
$ws.Range("H113").Value = 2253.1304
$ws.Range("I113").Value = 2283.8823
$ws.Range("J113").Value = 2166
$ws.Range("K113").Value = 2283.8823
$ws.Range("L113").Value = 2166
$ws.Range("M113").Value = -113.8823000000002
$ws.Range("N113").Value = -6506
$ws.Range("H139").Value = 84000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 84000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 84000
$ws.Range("N139").Value = -94280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10310.4
$ws.Range("I62").Value = 7850.6665
$ws.Range("J62").Value = 14000
$ws.Range("K62").Value = 7850.6665
$ws.Range("L62").Value = 14000
$ws.Range("M62").Value = -7226.6665
$ws.Range("N62").Value = -15248
$ws.Range("H65").Value = 10310.4
$ws.Range("I65").Value = 7850.6665
$ws.Range("J65").Value = 14000
$ws.Range("K65").Value = 39253.3325
$ws.Range("L65").Value = 70000
$ws.Range("M65").Value = -36133.3325
$ws.Range("N65").Value = -76240
$ws.Range("H74").Value = 16324.8
$ws.Range("I74").Value = 14575
$ws.Range("J74").Value = 17491.334
$ws.Range("K74").Value = 14575
$ws.Range("L74").Value = 17491.334
$ws.Range("M74").Value = -13639
$ws.Range("N74").Value = -19363.334
$ws.Range("H77").Value = 16324.8
$ws.Range("I77").Value = 14575
$ws.Range("J77").Value = 17491.334
$ws.Range("K77").Value = 43725
$ws.Range("L77").Value = 52474.00199999999
$ws.Range("M77").Value = -39045
$ws.Range("N77").Value = -61834.00199999999
$ws.Range("H81").Value = 4789.6
$ws.Range("I81").Value = 5155.1113
$ws.Range("J81").Value = 1500
$ws.Range("K81").Value = 10310.2226
$ws.Range("L81").Value = 3000
$ws.Range("M81").Value = -9249.222599999999
$ws.Range("N81").Value = -5122
$ws.Range("H84").Value = 4789.6
$ws.Range("I84").Value = 5155.1113
$ws.Range("J84").Value = 1500
$ws.Range("K84").Value = 51551.113
$ws.Range("L84").Value = 15000
$ws.Range("M84").Value = -46247.113
$ws.Range("N84").Value = -25608
$ws.Range("H108").Value = 100625.5
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 100625.5
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 100625.5
$ws.Range("N108").Value = -108305.5
$ws.Range("H122").Value = 3698.3057
$ws.Range("I122").Value = 2937.2917
$ws.Range("J122").Value = 5220.3335
$ws.Range("K122").Value = 8811.875100000001
$ws.Range("L122").Value = 15661.0005
$ws.Range("M122").Value = -6361.875100000001
$ws.Range("N122").Value = -20561.0005
